$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the same serial date value (45181) for every
# data row (rows 2-117). The update bumps that value by one day to 45182
# for all of those rows, leaving every other cell untouched.
for ($row = 2; $row -le 117; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45181) {
        $cell.Value = 45182
    }
}
